# Generate Report for Archive
#
# The localization status moves from "Ready for handoff" to "In Translation"
# for both locales tracked in this report, and the now-narrower status text
# means the "Status" columns can be shrunk accordingly.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Update the per-locale status cells on the Overview sheet (zh-cn -> column E,
# de-de -> column F) as well as the "Status" column (C) on each locale sheet.
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# Narrow the affected "Status" columns now that they no longer need to fit
# the longer "Ready for handoff" text.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
